$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 157
$ws.Range("J2").Value = 569
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 160
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 94
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 74
$ws.Range("T2").Value = 103
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 938
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 898
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 3
